$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4 (shifts rows 4..13 down to 5..14)
$ws.Rows(4).Insert()

# New row 4 is a copy of the surrounding records with updated date / price / origin values
$ws.Cells.Item(4, 1).Value = 10
$ws.Cells.Item(4, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(4, 3).Value = "La Araucanía"
$ws.Cells.Item(4, 4).Value = (Get-Date -Year 2023 -Month 12 -Day 14 -Hour 0 -Minute 0 -Second 0).Date
$ws.Cells.Item(4, 5).Value = 9
$ws.Cells.Item(4, 6).Value = "Fruta"
$ws.Cells.Item(4, 7).Value = 100101
$ws.Cells.Item(4, 8).Value = "Berries"
$ws.Cells.Item(4, 9).Value = 100101004
$ws.Cells.Item(4, 10).Value = "Frambuesa"
$ws.Cells.Item(4, 11).Value = "Sin especificar"
$ws.Cells.Item(4, 12).Value = "Primera"
$ws.Cells.Item(4, 13).Value = 40
$ws.Cells.Item(4, 14).Value = 7000
$ws.Cells.Item(4, 15).Value = 7000
$ws.Cells.Item(4, 16).Value = 7000
$ws.Cells.Item(4, 17).Value = '$/envase 1 kilo'
$ws.Cells.Item(4, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(4, 19).Value = 7000
$ws.Cells.Item(4, 20).Value = 1
